$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Dignity office admin removed": delete the entire row for
# Dignity Health / 3100 / Office Admin Provider Delegate / dignity.familypractice
# (original row 27). Everything below shifts up by one row; the now-unused
# shared string "dignity.familypractice" is dropped automatically by the
# recalculated shared-string table.
$ws.Rows.Item(27).Delete()

# The row delete above does not automatically renumber the worksheet's
# hyperlinks in this runtime, so rebuild the hyperlink collection to match
# the post-delete layout: the hyperlink that lived on the deleted row is
# gone, and every hyperlink that was below row 27 now lives one row higher.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D24"), "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=37053699&forPersonCustId=3200", "", "", "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=37053699&forPersonCustId=3200")
$ws.Hyperlinks.Add($ws.Range("D25"), "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=25530467&forPersonCustId=3100", "", "", "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=25530467&forPersonCustId=3100")
$ws.Hyperlinks.Add($ws.Range("D26"), "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=27542090&forPersonCustId=3100", "", "", "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=27542090&forPersonCustId=3100")

$ws.Hyperlinks.Add($ws.Range("D43"), "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=33993491&forPersonCustId=3500", "", "", "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=33993491&forPersonCustId=3500")
$ws.Hyperlinks.Add($ws.Range("D44"), "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=31837442&forPersonCustId=3500", "", "", "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=31837442&forPersonCustId=3500")
$ws.Hyperlinks.Add($ws.Range("D45"), "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=37034624&forPersonCustId=3500", "", "", "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=37034624&forPersonCustId=3500")

$ws.Hyperlinks.Add($ws.Range("D91"), "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=36963206&forPersonCustId=3800", "", "", "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=36963206&forPersonCustId=3800")
$ws.Hyperlinks.Add($ws.Range("D92"), "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=40707817&forPersonCustId=3800", "", "", "https://www.cozeva.com/new_settings?session=YXBwX2lkPXJlZ2lzdHJpZXMmY3VzdElkPTE1MDAmcGF5ZXJJZD0xNTAwJm9yZ0lkPTE1MDA&person_id=40707817&forPersonCustId=3800")

$ws.Hyperlinks.Add($ws.Range("D110"), "http://puchakraborty.cm/", "", "", "http://puchakraborty.cm/")

$ws.Hyperlinks.Add($ws.Range("D21"), "mailto:alex.turchinsky@cchphealthplan.com")

$ws.Hyperlinks.Add($ws.Range("A72"), "javascript:void(0);", "", "Optum Care Network - Idaho", "javascript:void(0);")
$ws.Hyperlinks.Add($ws.Range("A80"), "javascript:void(0);", "", "Optum Kansas City", "javascript:void(0);")

# Update the view to match the post-edit selection/scroll position.
$ws.Range("J24").Select()
$excel.ActiveWindow.ScrollRow = 4
